# Adjust position of plots in cheatsheet
#
# 1) Merge the two adjacent runs that together spelled out `"documents")`
#    into a single run with text `"documents")`.
# 2) Nudge the four plot pictures (bottom-right stack on slide 2) to their
#    new positions/sizes.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- 1) Merge the quote + "documents")" runs into a single run --------------
$sh = $s.Shapes.Item(6)
$tf = $sh.TextFrame
$tr = $tf.TextRange
$full = $tr.Text
$needle = '"documents")'
$idx = $full.IndexOf($needle)
if ($idx -ge 0) {
    $sub = $tr.Characters($idx + 1, $needle.Length)
    $sub.Text = $needle
}

# --- 2) Reposition / resize the four plot pictures --------------------------
# (point values below are the minimal doubles that round-trip through the
# float32-backed Left/Top/Width/Height properties to the exact target EMU)

# Picture 56
$pic = $s.Shapes.Item(13)
$pic.Left = 930.3065490722657
$pic.Top = 428.06526184082037

# Picture 57
$pic = $s.Shapes.Item(14)
$pic.Left = 942.0008850097657
$pic.Top = 325.8311004638672

# Picture 58
$pic = $s.Shapes.Item(15)
$pic.Left = 934.3443298339844
$pic.Width = 149.7940139770508
$pic.Height = 93.62125778198242

# Picture 59
$pic = $s.Shapes.Item(16)
$pic.Left = 930.3065490722657
$pic.Top = 646.0508728027344
